$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.191.90"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "2.359.30"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "244.75"
$ws.Range("E5").Value = "  +3.96%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "0.682"
$ws.Range("E6").Value = "  +4.76%  "
$ws.Range("D7").Value = "74.48"
$ws.Range("E7").Value = "  +3.38%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +24.71%  "
$ws.Range("E10").Value = "  +5.28%  "
$ws.Range("D11").Value = "32.13"
$ws.Range("E11").Value = "  +20.25%  "
$ws.Range("D12").Value = "7.51"
$ws.Range("E12").Value = "  +20.05%  "
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "2.711.73"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "16.93"
$ws.Range("E15").Value = "  +6.99%  "
$ws.Range("D16").Value = "0.916"
$ws.Range("E16").Value = "  +7.10%  "
$ws.Range("D17").Value = "2.357.88"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "44.415.59"
$ws.Range("E18").Value = "  +2.24%  "
$ws.Range("E19").Value = "  +4.10%  "
$ws.Range("D20").Value = "6.78"
$ws.Range("E20").Value = "  +6.12%  "
$ws.Range("D21").Value = "78.42"
$ws.Range("E21").Value = "  +5.22%  "
$ws.Range("D22").Value = "256.36"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "2.59"
$ws.Range("E23").Value = "  +4.95%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -5.08%  "
$ws.Range("D26").Value = "10.77"
$ws.Range("E26").Value = "  +7.41%  "
$ws.Range("E27").Value = "  +3.54%  "
$ws.Range("D28").Value = "22.58"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "1.61"
$ws.Range("E29").Value = "  +4.78%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'174.90"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("E31").Value = "  +4.09%  "
$ws.Range("E32").Value = "  +4.79%  "
$ws.Range("D33").Value = "5.43"
$ws.Range("E33").Value = "  +8.51%  "
$ws.Range("D34").Value = "'0.0760"
$ws.Range("E34").Value = "  +9.56%  "
$ws.Range("E35").Value = "  +5.68%  "
$ws.Range("E36").Value = "  +6.68%  "
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("E39").Value = "  +7.70%  "
$ws.Range("D40").Value = "19.48"
$ws.Range("E40").Value = "  +4.43%  "
$ws.Range("D41").Value = "9.01"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("E43").Value = "  +14.92%  "
$ws.Range("E44").Value = "  +3.24%  "
$ws.Range("D45").Value = "'2.50"
$ws.Range("E45").Value = "  +11.32%  "
$ws.Range("D46").Value = "0.0997"
$ws.Range("E46").Value = "  +4.74%  "
$ws.Range("D47").Value = "101.71"
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.455.15"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").Value = "0.000208"
$ws.Range("E51").Value = "  +5.13%  "
